$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New homework grade entered for student #2 (row 6), column I ---
# I6 needs the same "end of block" formatting (green fill + thick side
# border) that other cells in this column style already use (e.g. G15).
$ws.Range("G15").Copy() | Out-Null
$ws.Range("I6").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("I6").Value = 5

# --- Student #13 (row 17): grades filled in for C, D, E ---
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = 5

# --- Student #29 (row 33): grade filled in for F, new bonus grade H,
#     and the overall mark (column M) recorded as 3 ---
$ws.Range("F33").Value = 5

# H33 is a brand-new "bonus" cell; copy the bonus-column formatting
# (white fill + thick side border) from the existing K6 bonus cell.
$ws.Range("K6").Copy() | Out-Null
$ws.Range("H33").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("H33").Value = 5

$ws.Range("M33").Value = 3

# --- Selection moved to K6 ---
$ws.Range("K6").Select() | Out-Null

$excel.CutCopyMode = 0
